$d = $word.ActiveDocument

function InsertSourceBlock($paraIndex, $runs) {
    $p = $d.Paragraphs($paraIndex)
    $ip = $p.Range
    $ip.Collapse(1)
    $ip.InsertParagraphBefore()
    $newp = $d.Paragraphs($paraIndex)
    $newp.Style = "SourceCode"
    $pos = $newp.Range.Start
    foreach ($run in $runs) {
        $cur = $d.Range($pos, $pos)
        $cur.InsertAfter($run.Text)
        $endPos = $pos + $run.Text.Length
        if ($run.Style -ne '') {
            $styled = $d.Range($pos, $endPos)
            $styled.Style = $run.Style
        }
        $pos = $endPos
    }
}

function ReplaceParaLeadText($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $start = $p.Range.Start
    $r = $d.Range($start, $start + $oldText.Length)
    $r.Text = $newText
}

$runsSqrt = @(
    @{Text='x '; Style='NormalTok'},
    @{Text='<-'; Style='OtherTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='sample'; Style='FunctionTok'},
    @{Text='('; Style='NormalTok'},
    @{Text='2'; Style='DecValTok'},
    @{Text=':'; Style='SpecialCharTok'},
    @{Text='8'; Style='DecValTok'},
    @{Text=', '; Style='NormalTok'},
    @{Text='1'; Style='DecValTok'},
    @{Text=')'; Style='NormalTok'},
)

$runsOptsP = @(
    @{Text='opts_p '; Style='NormalTok'},
    @{Text='<-'; Style='OtherTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='c'; Style='FunctionTok'},
    @{Text='('; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='   '; Style='NormalTok'},
    @{Text='"the probability that the null hypothesis is true"'; Style='StringTok'},
    @{Text=','; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='   '; Style='NormalTok'},
    @{Text='answer ='; Style='AttributeTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='"the probability of the observed, or more extreme, data, under the assumption that the null-hypothesis is true"'; Style='StringTok'},
    @{Text=','; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='   '; Style='NormalTok'},
    @{Text='"the probability of making an error in your conclusion"'; Style='StringTok'},
    @{Text=([string][char]11); Style=''},
    @{Text=')'; Style='NormalTok'},
)

$runsOptsCi = @(
    @{Text='# use sample() to randomise the order'; Style='CommentTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='opts_ci '; Style='NormalTok'},
    @{Text='<-'; Style='OtherTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='sample'; Style='FunctionTok'},
    @{Text='('; Style='NormalTok'},
    @{Text='c'; Style='FunctionTok'},
    @{Text='('; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='  '; Style='NormalTok'},
    @{Text='answer ='; Style='AttributeTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='"if you repeated the process many times, 95% of intervals calculated in this way contain the true mean"'; Style='StringTok'},
    @{Text=','; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='  '; Style='NormalTok'},
    @{Text='"there is a 95% probability that the true mean lies within this range"'; Style='StringTok'},
    @{Text=','; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='  '; Style='NormalTok'},
    @{Text='"95% of the data fall within this range"'; Style='StringTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='))'; Style='NormalTok'},
)

$runsOpts = @(
    @{Text='opts '; Style='NormalTok'},
    @{Text='<-'; Style='OtherTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='c'; Style='FunctionTok'},
    @{Text='('; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='   '; Style='NormalTok'},
    @{Text='"the probability that the null hypothesis is true"'; Style='StringTok'},
    @{Text=','; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='   '; Style='NormalTok'},
    @{Text='answer ='; Style='AttributeTok'},
    @{Text=' '; Style='NormalTok'},
    @{Text='"the probability of the observed, or more extreme, data, under the assumption that the null-hypothesis is true"'; Style='StringTok'},
    @{Text=','; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text='   '; Style='NormalTok'},
    @{Text='"the probability of making an error in your conclusion"'; Style='StringTok'},
    @{Text=([string][char]11); Style=''},
    @{Text=')'; Style='NormalTok'},
    @{Text=([string][char]11); Style=''},
    @{Text=([string][char]11); Style=''},
    @{Text='cat'; Style='FunctionTok'},
    @{Text='('; Style='NormalTok'},
    @{Text='"What is a p-value?"'; Style='StringTok'},
    @{Text=', '; Style='NormalTok'},
    @{Text='longmcq'; Style='FunctionTok'},
    @{Text='(opts))'; Style='NormalTok'},
)

$runsCars = @(
    @{Text='with'; Style='FunctionTok'},
    @{Text='(cars, '; Style='NormalTok'},
    @{Text='plot'; Style='FunctionTok'},
    @{Text='(speed, dist))'; Style='NormalTok'},
)

# ===========================================================================
# Apply edits from bottom of the document to the top, so earlier paragraph
# indices are not invalidated by insertions that happen further down.
# ===========================================================================

# --- 8. with(cars, plot(speed, dist)) source block before the scatterplot
#        paragraph, and that paragraph's style BodyText -> FirstParagraph.
InsertSourceBlock 69 $runsCars
$pDraw = $d.Paragraphs(70)
$pDraw.Style = "FirstParagraph"

# --- 7. "What is a p-value?" (Checked sections copy): BodyText -> SourceCode
#        with the opts/cat() code, plus a new FirstParagraph paragraph after
#        it that still reads "What is a p-value?".
$p58 = $d.Paragraphs(58)
$p58.Style = "SourceCode"
$p58.Range.Text = ""
$pos = $p58.Range.Start
foreach ($run in $runsOpts) {
    $cur = $d.Range($pos, $pos)
    $cur.InsertAfter($run.Text)
    $endPos = $pos + $run.Text.Length
    if ($run.Style -ne '') {
        $styled = $d.Range($pos, $endPos)
        $styled.Style = $run.Style
    }
    $pos = $endPos
}
$afterCursor = $d.Range($p58.Range.End - 1, $p58.Range.End - 1)
$afterCursor.InsertParagraphAfter()
$pPval = $d.Paragraphs(59)
$pPval.Style = "FirstParagraph"
$pPval.Range.Text = "What is a p-value?"

# --- 6. Rotate the three confidence-interval MCQ option texts.
ReplaceParaLeadText 54 "if you repeated the process many times, 95% of intervals calculated in this way contain the true mean" "there is a 95% probability that the true mean lies within this range"
ReplaceParaLeadText 52 "95% of the data fall within this range" "if you repeated the process many times, 95% of intervals calculated in this way contain the true mean"
ReplaceParaLeadText 50 "there is a 95% probability that the true mean lies within this range" "95% of the data fall within this range"

# --- 5. opts_ci <- sample(c(...)) source block before the confidence
#        interval MCQ option list.
InsertSourceBlock 49 $runsOptsCi

# --- 4. opts_p <- c(...) source block before the p-value MCQ option list.
InsertSourceBlock 42 $runsOptsP

# --- 3 & 2. x <- sample(2:8, 1) source block before the "square root"
#        fill-in-the-blank, and 16 -> 25 in its text.
InsertSourceBlock 10 $runsSqrt
ReplaceParaLeadText 11 "The square root of 16 is: _" "The square root of 25 is: _"

Write-Output "done"
